$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values -----------------------------------------------
# Row 2 (earthquake)
$ws.Range("B2").Value = "deterministic"
$ws.Range("D2").Value = 0.001
$ws.Range("E2").Value = 0.15
$ws.Range("F2").Value = '{"landArea" : 2e-7}'
$ws.Range("G2").Value = '{"population" : -.1, "landArea" : -0.05}'

# Row 3 (hurricane) - Type changes from natural -> deterministic
$ws.Range("B3").Value = "deterministic"
$ws.Range("F3").Value = '{"landArea" : 1.5e-2}'
$ws.Range("G3").Value = '{"population" : -0.05, "timber" : -0.2, "housing" : -0.25}'

# Row 4 (drought) - Type changes from natural -> deterministic
$ws.Range("B4").Value = "deterministic"
$ws.Range("F4").Value = '{"food" : 1e-5}'
$ws.Range("G4").Value = '{"food" : -0.25, "landArea" : -0.05}'

# --- Remove the trailing blank row --------------------------------------
$ws.Rows.Item(5).Delete()

# --- Column widths -------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 19.59765625
$ws.Columns.Item(7).ColumnWidth = 50

# --- Selection -------------------------------------------------------------
$ws.Range("F2").Select()
